# Mise a jour des resultats du script - append new rows 435-448
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(435, '2026-02-04', 'eaux souterraines', 28, 1),
    @(436, '2026-02-04', 'ruissellement', 28, 3),
    @(437, '2026-02-04', 'ruissellement', 29, 2),
    @(438, '2026-02-04', 'ruissellement', 33, 2),
    @(439, '2026-02-04', 'zone tampon', 35, 1),
    @(440, '2026-02-04', 'bonnes pratiques', 36, 1),
    @(441, '2026-02-04', 'ruissellement', 36, 1),
    @(442, '2026-02-04', 'eaux souterraines', 37, 1),
    @(443, '2026-02-04', 'ruissellement', 37, 1),
    @(444, '2026-02-04', 'ruissellement', 41, 1),
    @(445, '2026-02-04', 'zone tampon', 46, 1),
    @(446, '2026-02-04', 'eaux souterraines', 47, 2),
    @(447, '2026-02-04', 'ruissellement', 47, 5),
    @(448, '2026-02-04', 'bonnes pratiques', 47, 1)
)

$firstRow = 435
$lastRow = 448

# Column A holds plain text dates (e.g. "2026-02-04"). Excel normally
# auto-detects this pattern and coerces it to a date serial number, so
# force a text number format before writing the values, then restore the
# default style afterwards so no visible formatting difference remains.
$dateRange = $ws.Range("A" + $firstRow + ":A" + $lastRow)
$dateRange.NumberFormat = "@"

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A" + $r).Value = $row[1]
    $ws.Range("B" + $r).Value = $row[2]
    $ws.Range("C" + $r).Value = $row[3]
    $ws.Range("D" + $r).Value = $row[4]
}

$dateRange.Style = "Normal"
